$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.474.10'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +4.63%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.590.52'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.27%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.61%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.73'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.68%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +1.06%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.63%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.87'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +8.17%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +1.25%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0602'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.43%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +2.12%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.816.80'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.24%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.600.67'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.75%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.08%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +2.05%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '28.471.73'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +4.85%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.09'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +2.99%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '233.31'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +7.81%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.52'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.35%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.90%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.72%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.37%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +2.04%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.37%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.53'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.92%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +1.66%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.33%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.107'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.94%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.65%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.49%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.31%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.08%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.30%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.418.51'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -2.18%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.29%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -5.95%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.64%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0167'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.36%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.54'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +8.26%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +1.85%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.51%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.79%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.68'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -2.90%  '
$ws.Range('B44').Value = 'WEMIXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.979'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -2.43%  '
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.83'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +6.11%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '64.67'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.34%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.727.97'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +1.35%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '87.80'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +1.92%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +2.58%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.07%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '39.60'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +16.80%  '
